# Update cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.264.62'
$ws.Range("E2").Value = '  +9.72%  '
$ws.Range("D3").Value = '3.224.02'
$ws.Range("E3").Value = '  +4.25%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''397.73'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").Value = '''111.20'
$ws.Range("E7").Value = '  +2.93%  '
$ws.Range("D9").Value = '''0.618'
$ws.Range("E9").Value = '  +5.96%  '
$ws.Range("D10").Value = '''39.20'
$ws.Range("E10").Value = '  +6.05%  '
$ws.Range("D11").Value = '''0.0927'
$ws.Range("E11").Value = '  +8.31%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '3.730.55'
$ws.Range("E13").Value = '  +4.21%  '
$ws.Range("E14").Value = '  +4.30%  '
$ws.Range("D15").Value = '''19.06'
$ws.Range("E15").Value = '  +3.21%  '
$ws.Range("D16").Value = '3.223.36'
$ws.Range("E16").Value = '  +3.89%  '
$ws.Range("E17").Value = '  +4.51%  '
$ws.Range("D18").Value = '''11.07'
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("D19").Value = '56.138.83'
$ws.Range("E19").Value = '  +9.30%  '
$ws.Range("D20").Value = '''3.35'
$ws.Range("E20").Value = '  +3.40%  '
$ws.Range("E21").Value = '  +7.52%  '
$ws.Range("D22").Value = '''13.03'
$ws.Range("E22").Value = '  +4.95%  '
$ws.Range("D23").Value = '''296.94'
$ws.Range("E23").Value = '  +11.78%  '
$ws.Range("D24").Value = '''76.20'
$ws.Range("E24").Value = '  +8.77%  '
$ws.Range("D25").Value = '''3.22'
$ws.Range("E25").Value = '  +1.77%  '
$ws.Range("D26").Value = '''8.15'
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").Value = '''28.02'
$ws.Range("E27").Value = '  +3.22%  '
$ws.Range("D28").Value = '''7.41'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("E29").Value = '  +4.16%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E31").Value = '  +4.06%  '
$ws.Range("D32").Value = '''11.15'
$ws.Range("E32").Value = '  +6.09%  '
$ws.Range("D33").Value = '''0.0492'
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("D34").Value = '''36.66'
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("E35").Value = '  +3.46%  '
$ws.Range("D36").Value = '''51.39'
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '''3.55'
$ws.Range("E37").Value = '  +4.40%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''3.09'
$ws.Range("E39").Value = '  +23.89%  '
$ws.Range("D40").Value = '''135.44'
$ws.Range("E40").Value = '  +3.84%  '
$ws.Range("D41").Value = '''17.45'
$ws.Range("E41").Value = '  +4.62%  '
$ws.Range("E42").Value = '  +3.61%  '
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").Value = '''0.119'
$ws.Range("E44").Value = '  +3.23%  '
$ws.Range("E45").Value = '  -2.86%  '
$ws.Range("D46").Value = '''22.27'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").Value = '''2.18'
$ws.Range("E47").Value = '  +52.93%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '''2.09'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.126.79'
$ws.Range("E49").Value = '  +2.59%  '
$ws.Range("D50").Value = '''2.43'
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("D51").Value = '''0.0360'
$ws.Range("E51").Value = '  +11.03%  '
